# A4 for jonathan. Samlet script til figurer.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Y" (sheet1.xml)
#   - shrink margins slightly on row 2/3, split the old "EH"/"Y" rows into
#     four rows (gas/oil added), and push the existing EH/Y rows down.
# ---------------------------------------------------------------------------
$wsY = $wb.Worksheets.Item("Y")

$wsY.Range("C2").Value = 0.49
$wsY.Range("C3").Value = 0.49

# Row 4 becomes the new Y_gas row
$wsY.Range("A4").Value = "Y_gas"
$wsY.Range("B4").Value = "Y"
$wsY.Range("C4").Value = 0.01

# Row 5 becomes the new Y_oil row
$wsY.Range("A5").Value = "Y_oil"
$wsY.Range("B5").Value = "Y"
$wsY.Range("C5").Value = 0.01

# Row 6 is the old row 4 (EH / out1 / 0.5), now shifted down
$wsY.Range("A6").Value = "EH"
$wsY.Range("B6").Value = "out1"
$wsY.Range("C6").Value = 0.5

# Row 7 is the old row 5 (Y / out1 / 0.5), now shifted down
$wsY.Range("A7").Value = "Y"
$wsY.Range("B7").Value = "out1"
$wsY.Range("C7").Value = 0.5

$wsY.Range("D11").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Q2P" (sheet2.xml)
#   - add Y_oil/oil and Y_gas/gas rows
# ---------------------------------------------------------------------------
$wsQ2P = $wb.Worksheets.Item("Q2P")

$wsQ2P.Range("A4").Value = "Y_oil"
$wsQ2P.Range("B4").Value = "oil"

$wsQ2P.Range("A5").Value = "Y_gas"
$wsQ2P.Range("B5").Value = "gas"

$wsQ2P.Range("B6").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "M" (sheet3.xml)
#   - update the oil abatement-cost figure
# ---------------------------------------------------------------------------
$wsM = $wb.Worksheets.Item("M")
$wsM.Range("E3").Value = 2.12

# Restore "M" as the active/visible sheet (it was active before this edit,
# and merely touching ranges on the other sheets above shifted focus away
# from it).
$wsM.Activate() | Out-Null
